$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted columns (B-E) stay as text, never auto-converted to numbers/dates
$ws.Range("B2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '26.746.67'
$ws.Range("E2").Value = '  -2.40%  '

# Row 3
$ws.Range("D3").Value = '1.796.84'
$ws.Range("E3").Value = '  -1.77%  '

# Row 4
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").Value = '308.21'
$ws.Range("E5").Value = '  -1.92%  '

# Row 6
$ws.Range("E6").Value = '  +0.09%  '

# Row 7
$ws.Range("D7").Value = '0.4573'
$ws.Range("E7").Value = '  +2.06%  '

# Row 8
$ws.Range("D8").Value = '0.3710'
$ws.Range("E8").Value = '  -1.80%  '

# Row 9
$ws.Range("D9").Value = '0.07246'
$ws.Range("E9").Value = '  -3.47%  '

# Row 10
$ws.Range("D10").Value = '0.8530'
$ws.Range("E10").Value = '  -5.15%  '

# Row 11
$ws.Range("D11").Value = '20.33'
$ws.Range("E11").Value = '  -3.42%  '

# Row 12
$ws.Range("D12").Value = '1.805.42'
$ws.Range("E12").Value = '  -1.24%  '

# Row 13
$ws.Range("D13").Value = '5.297'
$ws.Range("E13").Value = '  -2.16%  '

# Row 14
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").Value = '0.07037'
$ws.Range("E14").Value = '  -1.07%  '

# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '6.482'
$ws.Range("E15").Value = '  -4.32%  '

# Row 16
$ws.Range("D16").Value = '90.22'
$ws.Range("E16").Value = '  -4.48%  '

# Row 17
$ws.Range("E17").Value = '  +0.04%  '

# Row 18
$ws.Range("D18").Value = '0.000008624'
$ws.Range("E18").Value = '  -2.29%  '

# Row 19
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.19%  '

# Row 20
$ws.Range("D20").Value = '14.59'
$ws.Range("E20").Value = '  -4.33%  '

# Row 21
$ws.Range("D21").Value = '26.753.42'
$ws.Range("E21").Value = '  -2.44%  '

# Row 22
$ws.Range("D22").Value = '5.279'
$ws.Range("E22").Value = '  -0.30%  '

# Row 23
$ws.Range("D23").Value = '10.59'
$ws.Range("E23").Value = '  -3.36%  '

# Row 24
$ws.Range("D24").Value = '2.034.24'
$ws.Range("E24").Value = '  -0.96%  '

# Row 25
$ws.Range("D25").Value = '1.905'
$ws.Range("E25").Value = '  -5.10%  '

# Row 26
$ws.Range("D26").Value = '149.42'
$ws.Range("E26").Value = '  -1.44%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '18.17'
$ws.Range("E27").Value = '  -2.71%  '

# Row 28
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '2.143'
$ws.Range("E28").Value = '  -14.05%  '

# Row 29
$ws.Range("D29").Value = '5.200'
$ws.Range("E29").Value = '  -3.62%  '

# Row 30
$ws.Range("D30").Value = '114.20'
$ws.Range("E30").Value = '  -3.39%  '

# Row 31
$ws.Range("D31").Value = '0.08828'
$ws.Range("E31").Value = '  -0.34%  '

# Row 32
$ws.Range("D32").Value = '0.7516'
$ws.Range("E32").Value = '  -3.95%  '

# Row 33
$ws.Range("D33").Value = '1.156'
$ws.Range("E33").Value = '  -4.27%  '

# Row 34
$ws.Range("D34").Value = '4.431'
$ws.Range("E34").Value = '  -3.14%  '

# Row 35
$ws.Range("D35").Value = '2.880'
$ws.Range("E35").Value = '  -0.30%  '

# Row 36
$ws.Range("D36").Value = '1.000'
$ws.Range("E36").Value = '  +0.09%  '

# Row 37
$ws.Range("D37").Value = '1.114'
$ws.Range("E37").Value = '  +0.18%  '

# Row 38
$ws.Range("D38").Value = '0.01936'
$ws.Range("E38").Value = '  -2.83%  '

# Row 39
$ws.Range("D39").Value = '0.05209'
$ws.Range("E39").Value = '  -2.53%  '

# Row 40
$ws.Range("D40").Value = '2.902'
$ws.Range("E40").Value = '  +1.11%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '7.144'
$ws.Range("E41").Value = '  -3.40%  '

# Row 42
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '2.345'
$ws.Range("E42").Value = '  +2.82%  '

# Row 43
$ws.Range("D43").Value = '0.5209'
$ws.Range("E43").Value = '  -2.64%  '

# Row 44
$ws.Range("D44").Value = '0.1641'
$ws.Range("E44").Value = '  -5.26%  '

# Row 45
$ws.Range("D45").Value = '8.462'
$ws.Range("E45").Value = '  -4.19%  '

# Row 46
$ws.Range("D46").Value = '0.4985'
$ws.Range("E46").Value = '  -4.30%  '

# Row 47
$ws.Range("D47").Value = '10.23'
$ws.Range("E47").Value = '  -4.81%  '

# Row 48
$ws.Range("D48").Value = '103.99'
$ws.Range("E48").Value = '  -2.39%  '

# Row 49
$ws.Range("D49").Value = '1.000'
$ws.Range("E49").Value = '  +0.11%  '

# Row 50
$ws.Range("D50").Value = '1.643'
$ws.Range("E50").Value = '  -3.91%  '

# Row 51
$ws.Range("D51").Value = '0.06294'
$ws.Range("E51").Value = '  -1.38%  '
